$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts the existing FY data
# (old columns D:K, years 2011-2017) one column to the right (E:L),
# matching the workbook dimension growing from A5:K102 to A5:L102.
$ws.Columns("D:D").Insert()

# New column D inherits column C (text) formatting by default; copy the
# number/date formatting from column E (the old column D) into the new D
# column so dates keep the custom date format and figures keep the number format.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new column D with the latest (FY2018) period figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 9823000
$ws.Range("D9").Value = 2836000
$ws.Range("D10").Value = 6987000
$ws.Range("D12").Value = 1113000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 233000
$ws.Range("D15").Value = 599000
$ws.Range("D17").Value = 8317000
$ws.Range("D18").Value = 1506000
$ws.Range("D20").Value = 157000
$ws.Range("D21").Value = 2557000
$ws.Range("D22").Value = 241000
$ws.Range("D23").Value = 1422000
$ws.Range("D24").Value = -182000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1604000
$ws.Range("D27").Value = 1604000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 67000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -157000
$ws.Range("D33").Value = 1671000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1671000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 146000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 1608000
$ws.Range("D44").Value = 1166000
$ws.Range("D45").Value = 1083000
$ws.Range("D46").Value = 4003000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1782000
$ws.Range("D49").Value = 14283000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 931000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 20999000
$ws.Range("D57").Value = 349000
$ws.Range("D58").Value = 2253000
$ws.Range("D59").Value = 2658000
$ws.Range("D60").Value = 5260000
$ws.Range("D61").Value = 4803000
$ws.Range("D62").Value = 2210000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 12273000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -6953000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 8726000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 1671000
$ws.Range("D83").Value = 894000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 309000
$ws.Range("D91").Value = -316000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1921000
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 1432000
$ws.Range("D101").Value = -8000
$ws.Range("D102").Value = -188000

# Rows 9, 10 and 91 (Cost of Revenue, Gross Profit, Capital Expenditures) were
# restated for the prior periods as well (not a pure column shift) -- update
# columns E:J (2017..2012) with the revised figures; column K (2011, oldest
# year) is unchanged from the shift.
$ws.Range("E9").Value = 2616000
$ws.Range("F9").Value = 2469000
$ws.Range("G9").Value = 2212000
$ws.Range("H9").Value = 2297000
$ws.Range("I9").Value = 2314000
$ws.Range("J9").Value = 2494000

$ws.Range("E10").Value = 6432000
$ws.Range("F10").Value = 5917000
$ws.Range("G10").Value = 5265000
$ws.Range("H10").Value = 5083000
$ws.Range("I10").Value = 4829000
$ws.Range("J10").Value = 4755000

$ws.Range("E91").Value = -319000
$ws.Range("F91").Value = -376000
$ws.Range("G91").Value = -247000
$ws.Range("H91").Value = -259000
$ws.Range("I91").Value = -245000
$ws.Range("J91").Value = -226000

